$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("40-50","wanita","tidak","tidak","tidak","ya","ya","ya","tidak"),
    @("50-60","pria","tidak","tidak","ya","ya","tidak","ya","ya"),
    @("50-60","pria","ya","tidak","tidak","ya","ya","tidak","tidak"),
    @("50-60","wanita","ya","ya","tidak","ya","tidak","tidak","tidak"),
    @("50-60","pria","tidak","ya","ya","ya","tidak","ya","ya"),
    @("50-60","wanita","tidak","tidak","tidak","tidak","tidak","ya","tidak"),
    @("20-40","wanita","tidak","tidak","ya","ya","ya","ya","ya"),
    @("50-60","pria","ya","ya","ya","tidak","tidak","tidak","tidak"),
    @("40-50","pria","tidak","tidak","ya","tidak","ya","tidak","tidak"),
    @("50-60","wanita","ya","tidak","tidak","tidak","ya","tidak","ya"),
    @("50-60","wanita","ya","ya","tidak","ya","ya","ya","tidak"),
    @("50-60","wanita","ya","tidak","tidak","ya","tidak","ya","tidak"),
    @("40-50","wanita","tidak","tidak","ya","tidak","tidak","ya","tidak"),
    @("20-40","wanita","tidak","tidak","tidak","tidak","tidak","ya","ya"),
    @("50-60","pria","tidak","ya","ya","tidak","ya","ya","tidak")
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowValues[$j]
    }
}

$ws.Range("K7").Select()
